$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 6.312821333333335
$ws.Cells.Item(2, 8).Value = 18.938464
$ws.Cells.Item(2, 9).Value = 0.3104630857074662
$ws.Cells.Item(2, 10).Value = 0.3104630857074661
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 39.92590366666667
$ws.Cells.Item(2, 14).Value = 119.777711
$ws.Cells.Item(2, 15).Value = 0.8537340759835568
$ws.Cells.Item(2, 16).Value = 0.8537340759835568
$ws.Cells.Item(2, 17).Value = 252.045096419545
$ws.Cells.Item(2, 18).Value = 2268.405867775904
$ws.Cells.Item(2, 19).Value = 0.2650529156034674
$ws.Cells.Item(2, 20).Value = 0.2650529156034674

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 6.312821333333335
$ws.Cells.Item(3, 8).Value = 18.938464
$ws.Cells.Item(3, 9).Value = 0.3104630857074662
$ws.Cells.Item(3, 10).Value = 0.3104630857074661
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.145018
$ws.Cells.Item(3, 14).Value = 6.435054
$ws.Cells.Item(3, 15).Value = 0.04586683811810605
$ws.Cells.Item(3, 16).Value = 0.04586683811810605
$ws.Cells.Item(3, 17).Value = 13.541115390784
$ws.Cells.Item(3, 18).Value = 121.870038517056
$ws.Cells.Item(3, 19).Value = 0.01423996009379203
$ws.Cells.Item(3, 20).Value = 0.01423996009379203

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.312821333333335
$ws.Cells.Item(4, 8).Value = 18.938464
$ws.Cells.Item(4, 9).Value = 0.3104630857074662
$ws.Cells.Item(4, 10).Value = 0.3104630857074661
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.695284333333333
$ws.Cells.Item(4, 14).Value = 14.085853
$ws.Cells.Item(4, 15).Value = 0.1003990858983372
$ws.Cells.Item(4, 16).Value = 0.1003990858983372
$ws.Cells.Item(4, 17).Value = 29.64049110553245
$ws.Cells.Item(4, 18).Value = 266.7644199497921
$ws.Cells.Item(4, 19).Value = 0.03117021001020672
$ws.Cells.Item(4, 20).Value = 0.03117021001020672

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.354969666666667
$ws.Cells.Item(5, 8).Value = 28.064909
$ws.Cells.Item(5, 9).Value = 0.4600752335690602
$ws.Cells.Item(5, 10).Value = 0.4600752335690601
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 39.92590366666667
$ws.Cells.Item(5, 14).Value = 119.777711
$ws.Cells.Item(5, 15).Value = 0.8537340759835568
$ws.Cells.Item(5, 16).Value = 0.8537340759835568
$ws.Cells.Item(5, 17).Value = 373.5056177159221
$ws.Cells.Item(5, 18).Value = 3361.550559443299
$ws.Cells.Item(5, 19).Value = 0.3927819044140006
$ws.Cells.Item(5, 20).Value = 0.3927819044140006

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.354969666666667
$ws.Cells.Item(6, 8).Value = 28.064909
$ws.Cells.Item(6, 9).Value = 0.4600752335690602
$ws.Cells.Item(6, 10).Value = 0.4600752335690601
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.145018
$ws.Cells.Item(6, 14).Value = 6.435054
$ws.Cells.Item(6, 15).Value = 0.04586683811810605
$ws.Cells.Item(6, 16).Value = 0.04586683811810605
$ws.Cells.Item(6, 17).Value = 20.066578324454
$ws.Cells.Item(6, 18).Value = 180.599204920086
$ws.Cells.Item(6, 19).Value = 0.02110219626026191
$ws.Cells.Item(6, 20).Value = 0.02110219626026191

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.354969666666667
$ws.Cells.Item(7, 8).Value = 28.064909
$ws.Cells.Item(7, 9).Value = 0.4600752335690602
$ws.Cells.Item(7, 10).Value = 0.4600752335690601
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.695284333333333
$ws.Cells.Item(7, 14).Value = 14.085853
$ws.Cells.Item(7, 15).Value = 0.1003990858983372
$ws.Cells.Item(7, 16).Value = 0.1003990858983372
$ws.Cells.Item(7, 17).Value = 43.92424251470856
$ws.Cells.Item(7, 18).Value = 395.318182632377
$ws.Cells.Item(7, 19).Value = 0.04619113289479763
$ws.Cells.Item(7, 20).Value = 0.04619113289479762

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 4.665774
$ws.Cells.Item(8, 8).Value = 13.997322
$ws.Cells.Item(8, 9).Value = 0.2294616807234737
$ws.Cells.Item(8, 10).Value = 0.2294616807234737
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 39.92590366666667
$ws.Cells.Item(8, 14).Value = 119.777711
$ws.Cells.Item(8, 15).Value = 0.8537340759835568
$ws.Cells.Item(8, 16).Value = 0.8537340759835568
$ws.Cells.Item(8, 17).Value = 186.285243254438
$ws.Cells.Item(8, 18).Value = 1676.567189289942
$ws.Cells.Item(8, 19).Value = 0.1958992559660888
$ws.Cells.Item(8, 20).Value = 0.1958992559660887

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 4.665774
$ws.Cells.Item(9, 8).Value = 13.997322
$ws.Cells.Item(9, 9).Value = 0.2294616807234737
$ws.Cells.Item(9, 10).Value = 0.2294616807234737
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.145018
$ws.Cells.Item(9, 14).Value = 6.435054
$ws.Cells.Item(9, 15).Value = 0.04586683811810605
$ws.Cells.Item(9, 16).Value = 0.04586683811810605
$ws.Cells.Item(9, 17).Value = 10.008169213932
$ws.Cells.Item(9, 18).Value = 90.073522925388
$ws.Cells.Item(9, 19).Value = 0.0105246817640521
$ws.Cells.Item(9, 20).Value = 0.0105246817640521

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.665774
$ws.Cells.Item(10, 8).Value = 13.997322
$ws.Cells.Item(10, 9).Value = 0.2294616807234737
$ws.Cells.Item(10, 10).Value = 0.2294616807234737
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.695284333333333
$ws.Cells.Item(10, 14).Value = 14.085853
$ws.Cells.Item(10, 15).Value = 0.1003990858983372
$ws.Cells.Item(10, 16).Value = 0.1003990858983372
$ws.Cells.Item(10, 17).Value = 21.907135565074
$ws.Cells.Item(10, 18).Value = 197.164220085666
$ws.Cells.Item(10, 19).Value = 0.02303774299333287
$ws.Cells.Item(10, 20).Value = 0.02303774299333286
